$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.869.14'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '1.891.63'
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'0.7931"
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").Value = "'242.35"
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = "'0.3189"
$ws.Range("E8").Value = '  +2.13%  '
$ws.Range("D9").Value = "'25.71"
$ws.Range("E9").Value = '  -2.88%  '
$ws.Range("D10").Value = "'0.07087"
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = "'0.7747"
$ws.Range("E12").Value = '  +4.69%  '
$ws.Range("D13").Value = '1.965.89'
$ws.Range("E13").Value = '  +3.16%  '
$ws.Range("D14").Value = "'5.321"
$ws.Range("E14").Value = '  +2.78%  '
$ws.Range("D15").Value = "'92.37"
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").Value = '29.863.65'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = "'13.88"
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = "'5.938"
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").Value = "'244.20"
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("D20").Value = "'0.000007723"
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Value = "'8.132"
$ws.Range("E22").Value = '  +17.75%  '
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.140.07'
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = "'0.1629"
$ws.Range("E25").Value = '  +14.02%  '
$ws.Range("D26").Value = "'9.325"
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("D27").Value = "'164.82"
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("D28").Value = "'18.72"
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("D29").Value = "'2.070"
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").Value = "'1.374"
$ws.Range("E30").Value = '  +1.42%  '
$ws.Range("D31").Value = "'1.539"
$ws.Range("E31").Value = '  +1.69%  '
$ws.Range("D32").Value = "'4.437"
$ws.Range("E32").Value = '  +3.38%  '
$ws.Range("D33").Value = "'0.05658"
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("D34").Value = "'4.105"
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("D35").Value = "'1.270"
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = "'0.7374"
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").Value = "'2.710"
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").Value = "'0.01933"
$ws.Range("E39").Value = '  +0.34%  '
$ws.Range("D40").Value = "'2.782"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").Value = "'0.4445"
$ws.Range("E41").Value = '  +0.73%  '
$ws.Range("D42").Value = "'72.95"
$ws.Range("E42").Value = '  +1.16%  '
$ws.Range("D43").Value = "'5.854"
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("D44").Value = "'0.8435"
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.030.27'
$ws.Range("E46").Value = '  +5.12%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = "'1.884"
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = "'102.03"
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("D49").Value = "'9.936"
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("D50").Value = "'7.487"
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = "'2.960"
$ws.Range("E51").Value = '  +6.86%  '
